$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The chart source data only needs a single representative column (B) now;
# drop the extra per-car columns C through I.
$ws.Range("C1:I6").EntireColumn.Delete()

# Refresh the remaining summary values in column B.
$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 55.9

# B4 must stay a text value (like "43.82" before), not be reinterpreted as a number.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "24.23"
$ws.Range("B4").ClearFormats()

$ws.Range("B5").Value = 1508.5
$ws.Range("B6").Value = 10
